# Weekly data refresh: a new daily price observation (date serial 45106)
# is inserted as a new row ahead of the existing "Rabanito" series, pushing
# the rest of the historical rows down by one (row 419 -> 420, 420 -> 421,
# ..., 449 -> 450). The sheet's used-range / dimension grows from
# A1:R449 to A1:R450 automatically as part of the row insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at 419; Excel shifts rows 419..449 down to 420..450
# and keeps the sheet's dimension/used-range in sync.
$ws.Rows.Item(419).Insert()

# Populate the newly inserted row 419 with the new observation.
$ws.Cells.Item(419, 1).Value = 9
$ws.Cells.Item(419, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(419, 3).Value = "Metropolitana"
$ws.Cells.Item(419, 4).Value = 45106
$ws.Cells.Item(419, 5).Value = 13
$ws.Cells.Item(419, 6).Value = 300000001
$ws.Cells.Item(419, 7).Value = "Rabanito"
$ws.Cells.Item(419, 8).Value = "Sin especificar"
$ws.Cells.Item(419, 9).Value = "Primera"
$ws.Cells.Item(419, 10).Value = 7000
$ws.Cells.Item(419, 11).Value = 3000
$ws.Cells.Item(419, 12).Value = 3000
$ws.Cells.Item(419, 13).Value = 3000
$ws.Cells.Item(419, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(419, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(419, 16).Value = 30
$ws.Cells.Item(419, 17).Value = 100
$ws.Cells.Item(419, 18).Value = "Hortaliza"
